# edit.ps1 -- applies the two changes described by the diff:
#   1. Slide 5, shape 3 ("Google Shape;304;p17") text run "100" -> "3"
#   2. The presentation's theme colour scheme swaps from the "Momentum"
#      palette to the "Default" palette (the OOXML edit literally swaps
#      the contents of theme1.xml/theme2.xml; the reachable COM surface
#      for this is the active ThemeColorScheme on the SlideMaster, which
#      is backed by theme2.xml -- so we repaint its 12 colours to match
#      the "Default" palette that used to live in theme1.xml).

$p = $ppt.ActivePresentation

# --- 1. Text edit ---------------------------------------------------------
$s = $p.Slides.Item(5)
$s.Shapes.Item(3).TextFrame.TextRange.Text = "3"

# --- 2. Theme colour swap (Momentum -> Default) ---------------------------
$theme = $p.SlideMaster.Theme
$cs = $theme.ThemeColorScheme

# Index -> (role, new RGB as 0xBBGGRR long, matching the "Default" scheme)
$cs.Item(1).RGB  = 0            # dk1      000000
$cs.Item(2).RGB  = 16777215     # lt1      FFFFFF
$cs.Item(3).RGB  = 5800213      # dk2      158158
$cs.Item(4).RGB  = 15987699     # lt2      F3F3F3
$cs.Item(5).RGB  = 13077765     # accent1  058DC7
$cs.Item(6).RGB  = 3322960      # accent2  50B432
$cs.Item(7).RGB  = 1791725      # accent3  ED561B
$cs.Item(8).RGB  = 61421        # accent4  EDEF00
$cs.Item(9).RGB  = 15059748     # accent5  24CBE5
$cs.Item(10).RGB = 7529828      # accent6  64E572
$cs.Item(11).RGB = 13369378     # hlink    2200CC
$cs.Item(12).RGB = 9116245      # folHlink 551A8B
